{"js": "// Replace the date and all 2-digit x 2-digit multiplication problems/answers\n// with the updated values, matching the commit's diff exactly.\nconst replacements = [\n  [\"2026-02-17 Tuesday\", \"2026-02-18 Wednesday\"],\n  [\"43\u00d785=3655\", \"38\u00d727=1026\"],\n  [\"21\u00d794=1974\", \"70\u00d784=5880\"],\n  [\"56\u00d741=2296\", \"31\u00d720=620\"],\n  [\"34\u00d758=1972\", \"91\u00d736=3276\"],\n  [\"37\u00d745=1665\", \"60\u00d763=3780\"],\n  [\"34\u00d753=1802\", \"15\u00d766=990\"],\n  [\"32\u00d783=2656\", \"47\u00d748=2256\"],\n  [\"28\u00d793=2604\", \"94\u00d776=7144\"],\n  [\"63\u00d785=5355\", \"24\u00d782=1968\"],\n  [\"17\u00d715=255\", \"61\u00d749=2989\"],\n  [\"73\u00d774=5402\", \"66\u00d786=5676\"],\n  [\"59\u00d769=4071\", \"97\u00d768=6596\"],\n  [\"81\u00d793=7533\", \"13\u00d713=169\"],\n  [\"37\u00d711=407\", \"62\u00d727=1674\"],\n  [\"55\u00d738=2090\", \"97\u00d764=6208\"],\n  [\"90\u00d742=3780\", \"21\u00d787=1827\"],\n  [\"37\u00d744=1628\", \"14\u00d764=896\"],\n  [\"86\u00d785=7310\", \"32\u00d760=1920\"],\n  [\"50\u00d775=3750\", \"66\u00d775=4950\"],\n  [\"62\u00d783=5146\", \"28\u00d725=700\"],\n  [\"25\u00d749=1225\", \"77\u00d779=6083\"],\n  [\"83\u00d715=1245\", \"53\u00d731=1643\"],\n  [\"14\u00d727=378\", \"72\u00d760=4320\"],\n  [\"50\u00d773=3650\", \"98\u00d781=7938\"],\n  [\"99\u00d769=6831\", \"12\u00d795=1140\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and all 2-digit x 2-digit multiplication problems/answers\n# with the updated values, matching the commit's diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = '2026-02-17 Tuesday'; new = '2026-02-18 Wednesday'},\n    @{old = '43\u00d785=3655'; new = '38\u00d727=1026'},\n    @{old = '21\u00d794=1974'; new = '70\u00d784=5880'},\n    @{old = '56\u00d741=2296'; new = '31\u00d720=620'},\n    @{old = '34\u00d758=1972'; new = '91\u00d736=3276'},\n    @{old = '37\u00d745=1665'; new = '60\u00d763=3780'},\n    @{old = '34\u00d753=1802'; new = '15\u00d766=990'},\n    @{old = '32\u00d783=2656'; new = '47\u00d748=2256'},\n    @{old = '28\u00d793=2604'; new = '94\u00d776=7144'},\n    @{old = '63\u00d785=5355'; new = '24\u00d782=1968'},\n    @{old = '17\u00d715=255'; new = '61\u00d749=2989'},\n    @{old = '73\u00d774=5402'; new = '66\u00d786=5676'},\n    @{old = '59\u00d769=4071'; new = '97\u00d768=6596'},\n    @{old = '81\u00d793=7533'; new = '13\u00d713=169'},\n    @{old = '37\u00d711=407'; new = '62\u00d727=1674'},\n    @{old = '55\u00d738=2090'; new = '97\u00d764=6208'},\n    @{old = '90\u00d742=3780'; new = '21\u00d787=1827'},\n    @{old = '37\u00d744=1628'; new = '14\u00d764=896'},\n    @{old = '86\u00d785=7310'; new = '32\u00d760=1920'},\n    @{old = '50\u00d775=3750'; new = '66\u00d775=4950'},\n    @{old = '62\u00d783=5146'; new = '28\u00d725=700'},\n    @{old = '25\u00d749=1225'; new = '77\u00d779=6083'},\n    @{old = '83\u00d715=1245'; new = '53\u00d731=1643'},\n    @{old = '14\u00d727=378'; new = '72\u00d760=4320'},\n    @{old = '50\u00d773=3650'; new = '98\u00d781=7938'},\n    @{old = '99\u00d769=6831'; new = '12\u00d795=1140'}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
